# Update Mappings 22 Ontologies
# Adds a new "VIMMP_DEF" column (F) to the mapping sheet, with an
# empty-list placeholder value "[]" for each of the three data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell
$ws.Range("F1").Value = "VIMMP_DEF"

# New data cells (default/general formatting, same as columns C/E)
$ws.Range("F2").Value = "[]"
$ws.Range("F3").Value = "[]"
$ws.Range("F4").Value = "[]"

# Match the header formatting used by the other header cells (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
